$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain decimal numbers need to be forced to
# remain text (matching the original inlineStr/shared-string cell type) by
# temporarily switching the cell to Text format, assigning the value, then
# clearing the format again so no stray style survives.
$textCells = @(
    "D5",
    "D6",
    "D9",
    "D10",
    "D14",
    "D18",
    "D19",
    "D20",
    "D21",
    "D23",
    "D28",
    "D29",
    "D31",
    "D32",
    "D33",
    "D34",
    "D35",
    "D36",
    "D38",
    "D40",
    "D41",
    "D46",
    "D48",
    "D49",
    "D50",
    "D51"
)
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

# Apply every cell value from the diff.
$ws.Range("D2").Value = '58.076.39'
$ws.Range("E2").Value = '  +1.95%  '
$ws.Range("D3").Value = '3.076.09'
$ws.Range("E3").Value = '  +0.87%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '516.23'
$ws.Range("E5").Value = '  +0.81%  '
$ws.Range("D6").Value = '141.50'
$ws.Range("E6").Value = '  +1.53%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("E8").Value = '  +0.73%  '
$ws.Range("D9").Value = '7.33'
$ws.Range("E9").Value = '  +2.38%  '
$ws.Range("D10").Value = '0.108'
$ws.Range("E10").Value = '  +0.23%  '
$ws.Range("E11").Value = '  +3.25%  '
$ws.Range("D12").Value = '3.603.13'
$ws.Range("E12").Value = '  +0.80%  '
$ws.Range("E13").Value = '  +1.37%  '
$ws.Range("D14").Value = '26.71'
$ws.Range("E14").Value = '  +6.46%  '
$ws.Range("E15").Value = '  +0.55%  '
$ws.Range("D16").Value = '58.065.19'
$ws.Range("E16").Value = '  +1.84%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.093.63'
$ws.Range("E17").Value = '  +1.35%  '
$ws.Range("B18").Value = 'Polkadot'
$ws.Range("C18").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D18").Value = '6.17'
$ws.Range("E18").Value = '  +4.68%  '
$ws.Range("D19").Value = '12.85'
$ws.Range("E19").Value = '  -1.16%  '
$ws.Range("D20").Value = '8.11'
$ws.Range("E20").Value = '  +0.29%  '
$ws.Range("D21").Value = '333.45'
$ws.Range("E21").Value = '  -0.09%  '
$ws.Range("E22").Value = '  -0.01%  '
$ws.Range("D23").Value = '0.502'
$ws.Range("E23").Value = '  +0.19%  '
$ws.Range("E24").Value = '  +0.06%  '
$ws.Range("E25").Value = '  +3.07%  '
$ws.Range("E26").Value = '  -0.13%  '
$ws.Range("E27").Value = '  -2.79%  '
$ws.Range("D28").Value = '6.45'
$ws.Range("E28").Value = '  +1.72%  '
$ws.Range("D29").Value = '7.21'
$ws.Range("E29").Value = '  +4.63%  '
$ws.Range("E30").Value = '  +0.44%  '
$ws.Range("D31").Value = '1.21'
$ws.Range("E31").Value = '  +3.37%  '
$ws.Range("D32").Value = '20.81'
$ws.Range("E32").Value = '  +0.83%  '
$ws.Range("D33").Value = '155.03'
$ws.Range("E33").Value = '  +0.27%  '
$ws.Range("D34").Value = '4.55'
$ws.Range("E34").Value = '  +1.39%  '
$ws.Range("D35").Value = '27.36'
$ws.Range("E35").Value = '  +3.61%  '
$ws.Range("D36").Value = '6.02'
$ws.Range("E36").Value = '  +3.55%  '
$ws.Range("E37").Value = '  +4.62%  '
$ws.Range("D38").Value = '0.0675'
$ws.Range("E38").Value = '  +1.46%  '
$ws.Range("D39").Value = '3.116.46'
$ws.Range("E39").Value = '  +0.94%  '
$ws.Range("D40").Value = '3.90'
$ws.Range("E40").Value = '  +2.63%  '
$ws.Range("D41").Value = '36.59'
$ws.Range("E41").Value = '  -0.60%  '
$ws.Range("E42").Value = '  +0.01%  '
$ws.Range("E43").Value = '  -1.54%  '
$ws.Range("D44").Value = '2.279.79'
$ws.Range("E44").Value = '  +2.21%  '
$ws.Range("E45").Value = '  +1.41%  '
$ws.Range("D46").Value = '1.39'
$ws.Range("E46").Value = '  +2.34%  '
$ws.Range("E47").Value = '  +4.19%  '
$ws.Range("D48").Value = '0.943'
$ws.Range("E48").Value = '  +1.35%  '
$ws.Range("D49").Value = '5.95'
$ws.Range("E49").Value = '  +2.40%  '
$ws.Range("D50").Value = '0.734'
$ws.Range("E50").Value = '  +7.71%  '
$ws.Range("D51").Value = '257.89'
$ws.Range("E51").Value = '  +11.19%  '

# Restore the temporarily-applied Text format so the cells end up with no
# explicit style again, just like in the original workbook.
foreach ($c in $textCells) {
    $ws.Range($c).ClearFormats()
}
